# Scheduled runner update: refresh Sheets market-price snapshot values
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns H:N)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 42944
$ws.Range("I64").Value = 93600
$ws.Range("J64").Value = 3142.8572
$ws.Range("K64").Value = 93600
$ws.Range("L64").Value = 3142.8572
$ws.Range("M64").Value = -93352
$ws.Range("N64").Value = -3638.8572
$ws.Range("H67").Value = 42944
$ws.Range("I67").Value = 93600
$ws.Range("J67").Value = 3142.8572
$ws.Range("K67").Value = 93600
$ws.Range("L67").Value = 3142.8572
$ws.Range("M67").Value = -92742
$ws.Range("N67").Value = -4858.8572
$ws.Range("H76").Value = 5059.6665
$ws.Range("I76").Value = 4695.8887
$ws.Range("J76").Value = 5332.5
$ws.Range("K76").Value = 4695.8887
$ws.Range("L76").Value = 5332.5
$ws.Range("M76").Value = -4380.8887
$ws.Range("N76").Value = -5962.5
$ws.Range("H79").Value = 5059.6665
$ws.Range("I79").Value = 4695.8887
$ws.Range("J79").Value = 5332.5
$ws.Range("K79").Value = 4695.8887
$ws.Range("L79").Value = 5332.5
$ws.Range("M79").Value = -3603.8887
$ws.Range("N79").Value = -7516.5
$ws.Range("H106").Value = 3873.7778
$ws.Range("I106").Value = 3858
$ws.Range("K106").Value = 3858
$ws.Range("M106").Value = -3227
$ws.Range("H129").Value = 6732.8823
$ws.Range("I129").Value = 14793.143
$ws.Range("J129").Value = 1090.7
$ws.Range("K129").Value = 44379.429
$ws.Range("L129").Value = 3272.1
$ws.Range("M129").Value = -39379.429
$ws.Range("N129").Value = -13272.1
$ws.Range("H138").Value = 3349.0132
$ws.Range("I138").Value = 2055.04
$ws.Range("J138").Value = 3983.3137
$ws.Range("K138").Value = 6165.12
$ws.Range("L138").Value = 11949.9411
$ws.Range("M138").Value = -1025.12
$ws.Range("N138").Value = -22229.9411

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35685.277
$ws.Range("I32").Value = 15090.552
$ws.Range("J32").Value = 189001.56
$ws.Range("K32").Value = 15090.552
$ws.Range("L32").Value = 189001.56
$ws.Range("M32").Value = -14803.552
$ws.Range("N32").Value = -189575.56
$ws.Range("H122").Value = 2495.6
$ws.Range("I122").Value = 2196.1667
$ws.Range("K122").Value = 6588.500100000001
$ws.Range("M122").Value = -4138.500100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31166.686
$ws.Range("I31").Value = 1437.5186
$ws.Range("K31").Value = 1437.5186
$ws.Range("M31").Value = -1142.5186
$ws.Range("H34").Value = 31166.686
$ws.Range("I34").Value = 1437.5186
$ws.Range("K34").Value = 1437.5186
$ws.Range("M34").Value = -1235.5186
$ws.Range("H41").Value = 11176
$ws.Range("I41").Value = 5750
$ws.Range("J41").Value = 14793.333
$ws.Range("K41").Value = 5750
$ws.Range("L41").Value = 14793.333
$ws.Range("M41").Value = -5322
$ws.Range("N41").Value = -15649.333
$ws.Range("H50").Value = 9435.429
$ws.Range("J50").Value = 9435.429
$ws.Range("L50").Value = 9435.429
$ws.Range("N50").Value = -10685.429
$ws.Range("H51").Value = 7917.8
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7917.8
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7917.8
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9389.799999999999
$ws.Range("H60").Value = 15012.6
$ws.Range("J60").Value = 15012.6
$ws.Range("L60").Value = 15012.6
$ws.Range("N60").Value = -16034.6
$ws.Range("H61").Value = 7917.8
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 7917.8
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 7917.8
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -8613.799999999999
$ws.Range("H62").Value = 2543.75
$ws.Range("I62").Value = 2425
$ws.Range("J62").Value = 2662.5
$ws.Range("K62").Value = 2425
$ws.Range("L62").Value = 2662.5
$ws.Range("M62").Value = -1801
$ws.Range("N62").Value = -3910.5
$ws.Range("H65").Value = 2543.75
$ws.Range("I65").Value = 2425
$ws.Range("J65").Value = 2662.5
$ws.Range("K65").Value = 12125
$ws.Range("L65").Value = 13312.5
$ws.Range("M65").Value = -9005
$ws.Range("N65").Value = -19552.5
$ws.Range("H68").Value = 17726.777
$ws.Range("J68").Value = 17726.777
$ws.Range("L68").Value = 17726.777
$ws.Range("N68").Value = -19224.777
$ws.Range("H71").Value = 17726.777
$ws.Range("J71").Value = 17726.777
$ws.Range("L71").Value = 53180.33099999999
$ws.Range("N71").Value = -60668.33099999999
$ws.Range("H74").Value = 40742.332
$ws.Range("J74").Value = 40742.332
$ws.Range("L74").Value = 40742.332
$ws.Range("N74").Value = -42490.332
$ws.Range("H77").Value = 40742.332
$ws.Range("J77").Value = 40742.332
$ws.Range("L77").Value = 122226.996
$ws.Range("N77").Value = -130962.996

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 342.22223
$ws.Range("I18").Value = 277.5
$ws.Range("K18").Value = 832.5
$ws.Range("M18").Value = -663.5
$ws.Range("H134").Value = 4397.1816
$ws.Range("I134").Value = 2603.4285
$ws.Range("J134").Value = 7536.25
$ws.Range("K134").Value = 7810.2855
$ws.Range("L134").Value = 22608.75
$ws.Range("M134").Value = -2740.2855
$ws.Range("N134").Value = -32748.75
$ws.Range("H139").Value = 2441.818
$ws.Range("I139").Value = 1401.4286
$ws.Range("J139").Value = 2927.3333
$ws.Range("K139").Value = 4204.2858
$ws.Range("L139").Value = 8781.999899999999
$ws.Range("M139").Value = 935.7142000000003
$ws.Range("N139").Value = -19061.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 137367.4
$ws.Range("I70").Value = 225912.33
$ws.Range("J70").Value = 4550
$ws.Range("K70").Value = 225912.33
$ws.Range("L70").Value = 4550
$ws.Range("M70").Value = -225642.33
$ws.Range("N70").Value = -5090
$ws.Range("H73").Value = 137367.4
$ws.Range("I73").Value = 225912.33
$ws.Range("J73").Value = 4550
$ws.Range("K73").Value = 225912.33
$ws.Range("L73").Value = 4550
$ws.Range("M73").Value = -224976.33
$ws.Range("N73").Value = -6422
$ws.Range("H80").Value = 167037740
$ws.Range("J80").Value = 10750
$ws.Range("L80").Value = 10750
$ws.Range("N80").Value = -12746
$ws.Range("H83").Value = 167037740
$ws.Range("J83").Value = 10750
$ws.Range("L83").Value = 53750
$ws.Range("N83").Value = -63734
$ws.Range("H126").Value = 4943.5713
$ws.Range("I126").Value = 4918.2
$ws.Range("J126").Value = 5007
$ws.Range("K126").Value = 14754.6
$ws.Range("L126").Value = 15021
$ws.Range("M126").Value = -12284.6
$ws.Range("N126").Value = -19961

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3418.238
$ws.Range("I93").Value = 3428.2942
$ws.Range("J93").Value = 3375.5
$ws.Range("K93").Value = 3428.2942
$ws.Range("L93").Value = 3375.5
$ws.Range("M93").Value = -2180.2942
$ws.Range("N93").Value = -5871.5
$ws.Range("H122").Value = 4547.933
$ws.Range("I122").Value = 3662.182
$ws.Range("K122").Value = 10986.546
$ws.Range("M122").Value = -8536.545999999998
$ws.Range("H140").Value = 59261
$ws.Range("J140").Value = 59261
$ws.Range("L140").Value = 59261
$ws.Range("N140").Value = -69621

Write-Host "Edits applied successfully."
